$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet (fund-holdings detail) right after "2021-Q4"
#    and right before "总计", by duplicating the "2021-Q4" sheet (so headers,
#    column-A styling and borders all match the existing pattern exactly)
#    and then overwriting its data with the 2022-Q1 figures.
# ---------------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$newSheet = $wb.Worksheets.Item("2021-Q4 (2)")
$newSheet.Name = "2022-Q1"

# The template sheet has 21 data rows (rows 2-22); 2022-Q1 only needs 12
# (rows 2-13), so drop the extra rows entirely.
$newSheet.Range("A14:H22").Clear()

$fundRows = @(
    @("470006", "汇添富医药保健混合A",             "49.54", "88.92", "3.58", "1.7735", 10),
    @("009664", "汇添富医疗积极成长一年持有期混合A", "38.71", "67.08", "3.75", "1.4516", 4),
    @("001725", "汇添富中国高端制造股票",             "18.26", "90.98", "5.30", "0.9678", 5),
    @("015115", "汇添富中国高端制造股票D",            "18.26", "90.98", "5.30", "0.9678", 5),
    @("008415", "国泰大制造两年持有期混合",           "23.19", "92.05", "2.69", "0.6238", 9),
    @("005823", "泰康颐享混合A",                     "14.39", "20.19", "1.54", "0.2216", 5),
    @("009665", "汇添富医疗积极成长一年持有期混合C", "4.37",  "67.08", "3.75", "0.1639", 4),
    @("001907", "国投瑞银境煊灵活配置混合A",         "2.61",  "90.44", "4.58", "0.1195", 8),
    @("001908", "国投瑞银境煊灵活配置混合C",         "1.75",  "90.44", "4.58", "0.0802", 8),
    @("005824", "泰康颐享混合C",                     "2.82",  "20.19", "1.54", "0.0434", 5),
    @("015309", "国投瑞银境煊灵活配置混合E",         "0.33",  "90.44", "4.58", "0.0151", 8),
    @("960015", "汇添富医药保健混合O",               "0.00",  "88.92", "3.58", $null,    10)
)

for ($i = 0; $i -lt $fundRows.Length; $i++) {
    $r = $i + 2
    $row = $fundRows[$i]

    $newSheet.Cells.Item($r, 1).Value = $i

    # Fund code / name / scale / stock-position / position-ratio are stored
    # as text in the source data (e.g. "470006", "49.54") - force text so
    # numeric-looking strings are not reinterpreted as numbers.
    $newSheet.Cells.Item($r, 2).NumberFormat = "@"
    $newSheet.Cells.Item($r, 2).Value = $row[0]
    $newSheet.Cells.Item($r, 3).NumberFormat = "@"
    $newSheet.Cells.Item($r, 3).Value = $row[1]
    $newSheet.Cells.Item($r, 4).NumberFormat = "@"
    $newSheet.Cells.Item($r, 4).Value = $row[2]
    $newSheet.Cells.Item($r, 5).NumberFormat = "@"
    $newSheet.Cells.Item($r, 5).Value = $row[3]
    $newSheet.Cells.Item($r, 6).NumberFormat = "@"
    $newSheet.Cells.Item($r, 6).Value = $row[4]

    if ($row[5] -eq $null) {
        # Last row's held-market-value is a genuine number (0), not text.
        $newSheet.Cells.Item($r, 7).Value = 0
    } else {
        $newSheet.Cells.Item($r, 7).NumberFormat = "@"
        $newSheet.Cells.Item($r, 7).Value = $row[5]
    }

    $newSheet.Cells.Item($r, 8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (summary) sheet: add a new 2022-Q1 row at the top of
#    the data and keep the rest (re-indexing column A as Excel would).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Extend column-A's bordered/bold style down to the new row 7 before writing
# values into it (row 7 did not exist before, so it has no style yet).
$total.Range("A2").Copy()
$total.Range("A7").PasteSpecial(-4122)

$summaryRows = @(
    @(0, "2022-Q1", 12, 6.43),
    @(1, "2021-Q4", 21, 9.5),
    @(2, "2021-Q3", 22, 10.88),
    @(3, "2021-Q2", 29, 14.98),
    @(4, "2021-Q1", 28, 12.34),
    @(5, "2020-Q4", 21, 10.93)
)

for ($i = 0; $i -lt $summaryRows.Length; $i++) {
    $r = $i + 2
    $row = $summaryRows[$i]
    $total.Cells.Item($r, 1).Value = $row[0]
    $total.Cells.Item($r, 2).Value = $row[1]
    $total.Cells.Item($r, 3).Value = $row[2]
    $total.Cells.Item($r, 4).Value = $row[3]
}

Write-Host "2022-Q1 sheet and 总计 summary updated"
